$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Widen column E to fit the new "保留位数后..." labels
$ws.Columns.Item(5).ColumnWidth = 18.33203125

# Move the footnote "注：此处会判断数据个数大于等于6时自动近似为t=1" from D15 up to D11
$ws.Range("D15").Value = $null
$ws.Range("D11").Value = "注：此处会判断数据个数大于等于6时自动近似为t=1"

# New row 14: headers for the rounded output cells
$ws.Range("E14").Value = "保留位数后(与μD对齐位数)"
$ws.Range("J14").Value = "保留位数后(最高次有效数字为1,2保留两位,反之保留一位)"

# New row 15: formulas computing the rounded/aligned output values
$ws.Range("E15").Formula = "=IF(J13*10^INT(-LOG(ABS(J13)))<0.3,FIXED(E13,1-INT(LOG(J13)),1),FIXED(E13,-INT(LOG(J13)),1))"
$ws.Range("J15").Formula = "=IF(J13*10^INT(-LOG(ABS(J13)))<0.3,FIXED(J13,1-INT(LOG(J13)),1),FIXED(J13,-INT(LOG(J13)),1))"

# Restore the active cell selection used before saving
$ws.Range("I6").Select()
